$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 0. Remove the existing "总计" sheet up front. All the downstream
#    sheet-id bookkeeping (sheetId = max(existing sheetIds) + 1) works
#    out so that the freshly inserted "2022-Q1" sheet reuses the id
#    that "总计" used to have, and the rebuilt "总计" sheet gets the
#    next id after that -- matching how Excel numbers sheets that are
#    inserted after a deletion.
# ------------------------------------------------------------------
$wsTotalOld = $wb.Worksheets.Item("总计")
$wsTotalOld.Delete()

# ------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet by copying the existing "2021-Q4"
#    sheet (rebuilding from scratch would lose the header/index-column
#    styling), then overwrite it with the 2022-Q1 fund data.
# ------------------------------------------------------------------
$ws2021Q4 = $wb.Worksheets.Item("2021-Q4")
$ws2021Q4.Copy($null, $ws2021Q4)
$wsNew = $wb.Worksheets.Item($ws2021Q4.Index + 1)
$wsNew.Name = "2022-Q1"

$fundRows = @(
    @('002560','诺安和鑫灵活配置混合','33.85','93.22','9.23','3.1244',3),
    @('010680','华夏新兴成长股票A','55.80','87.37','5.23','2.9183',5),
    @('001071','华安媒体互联网混合','51.61','92.88','4.62','2.3844',6),
    @('001694','华安沪港深外延增长灵活配置混合','43.58','92.63','4.67','2.0352',7),
    @('010305','华夏创新驱动混合A','28.56','88.57','5.49','1.5679',5),
    @('310358','申万菱信新经济混合','41.92','77.32','3.64','1.5259',8),
    @('006879','华安智能生活混合','33.48','92.78','4.46','1.4932',6),
    @('007872','金信稳健策略灵活配置混合','25.57','93.73','5.31','1.3578',8),
    @('161914','万家创业板2年定期开放混合A','15.74','95.20','6.20','0.9759',9),
    @('506002','易方达科创板两年定期开放混合型证券投资基金','29.21','86.29','3.19','0.9318',9),
    @('001513','易方达信息产业混合','32.50','92.37','2.81','0.9132',6),
    @('010013','易方达信息行业精选股票','25.38','88.64','3.54','0.8985',3),
    @('160325','华夏创业板两年定期开放混合','27.39','90.77','2.63','0.7204',6),
    @('010286','海富通成长价值混合A','19.35','92.83','3.68','0.7121',8),
    @('001404','招商移动互联网产业股票','13.45','90.96','4.99','0.6712',7),
    @('007460','华安成长创新混合','13.21','91.10','4.64','0.6129',6),
    @('002229','华夏经济转型股票','11.71','86.14','5.11','0.5984',4),
    @('519957','长信睿进灵活配置混合A','15.70','44.52','3.19','0.5008',3),
    @('519956','长信睿进灵活配置混合C','15.69','44.52','3.19','0.5005',3),
    @('519011','海富通精选混合','10.25','76.05','3.49','0.3577',9),
    @('001042','华夏领先股票','13.72','93.46','2.56','0.3512',6),
    @('010681','华夏新兴成长股票C','5.29','87.37','5.23','0.2767',5),
    @('006868','华夏科技成长股票','5.23','87.33','5.11','0.2673',5),
    @('008655','招商科技创新混合A','4.72','90.71','5.12','0.2417',4),
    @('013840','银华集成电路混合A','8.32','71.47','2.83','0.2355',10),
    @('012210','申万菱信智能汽车股票型证券投资基金A','4.76','82.52','4.27','0.2033',8),
    @('011669','长信优质企业混合型证券投资基金A','4.24','88.48','4.04','0.1713',9),
    @('160425','华安创业板两年定期开放混合','5.11','96.75','2.90','0.1482',8),
    @('161915','万家创业板2年定期开放混合C','2.36','95.20','6.20','0.1463',9),
    @('519015','海富通精选贰号混合','3.73','76.65','3.55','0.1324',9),
    @('002256','金信行业优选灵活配置混合','2.43','93.89','5.32','0.1293',8),
    @('013634','申万菱信双利混合A','7.83','22.26','1.63','0.1276',6),
    @('010306','华夏创新驱动混合C','2.15','88.57','5.49','0.1180',5),
    @('004314','前海开源沪港深新硬件主题灵活配置混合A','1.67','90.05','5.62','0.0939',3),
    @('008656','招商科技创新混合C','1.57','90.71','5.12','0.0804',4),
    @('010287','海富通成长价值混合C','1.81','92.83','3.68','0.0666',8),
    @('012211','申万菱信智能汽车股票型证券投资基金C','1.40','82.52','4.27','0.0598',8),
    @('006502','财通集成电路产业股票A','1.29','79.76','4.58','0.0591',7),
    @('004315','前海开源沪港深新硬件主题灵活配置混合C','1.00','90.05','5.62','0.0562',3),
    @('004558','汇安丰裕灵活配置混合A','0.99','83.41','2.94','0.0291',9),
    @('011670','长信优质企业混合型证券投资基金C','0.69','88.48','4.04','0.0279',9),
    @('006503','财通集成电路产业股票C','0.46','79.76','4.58','0.0211',7),
    @('013841','银华集成电路混合C','0.72','71.47','2.83','0.0204',10),
    @('013635','申万菱信双利混合C','0.75','22.26','1.63','0.0122',6),
    @('519935','长信创新驱动股票','0.21','90.78','3.79','0.0080',5),
    @('013903','国泰君安信息行业混合','0.25','84.06','2.84','0.0071',10),
    @('001978','泰信互联网+主题灵活配置混合','0.06','92.34','3.08','0.0018',6),
    @('004559','汇安丰裕灵活配置混合C','0.01','83.41','2.94','0.0003',9)
)

$row = 2
foreach ($r in $fundRows) {
    $wsNew.Range("B$row").Value = "'" + $r[0]
    $wsNew.Range("C$row").Value = "'" + $r[1]
    $wsNew.Range("D$row").Value = "'" + $r[2]
    $wsNew.Range("E$row").Value = "'" + $r[3]
    $wsNew.Range("F$row").Value = "'" + $r[4]
    $wsNew.Range("G$row").Value = "'" + $r[5]
    $wsNew.Range("H$row").Value = $r[6]
    $row = $row + 1
}

# The copied template ("2021-Q4") had 88 data rows (rows 2-89); the new
# sheet only needs 48 (rows 2-49), so clear out the leftover rows.
$wsNew.Range("A50:H89").Clear()

# ------------------------------------------------------------------
# 2. Rebuild the "总计" (summary) sheet. Copy "2021-Q4" again (for its
#    header/index-column styling), trim it down to 4 columns x 6 rows,
#    and fill in the summary data with the new leading 2022-Q1 row.
# ------------------------------------------------------------------
$wsNew.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTotalNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTotalNew.Name = "总计"

$wsTotalNew.Range("E1:H89").Clear()
$wsTotalNew.Range("A7:D89").Clear()

$wsTotalNew.Range("B1").Value = "日期"
$wsTotalNew.Range("C1").Value = "持有数量(只)"
$wsTotalNew.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @('2022-Q1', 48, 27.89),
    @('2021-Q4', 88, 49.91),
    @('2021-Q3', 55, 18.59),
    @('2021-Q2', 12, 6.22),
    @('2020-Q4', 1, 0.01)
)

$row = 2
$idx = 0
foreach ($r in $summaryRows) {
    $wsTotalNew.Range("A$row").Value = $idx
    $wsTotalNew.Range("B$row").Value = "'" + $r[0]
    $wsTotalNew.Range("C$row").Value = $r[1]
    $wsTotalNew.Range("D$row").Value = $r[2]
    $row = $row + 1
    $idx = $idx + 1
}
